$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would be auto-parsed by Excel as a Number if assigned
# directly (e.g. "226.23"). The source sheet keeps these as plain text cells
# (t="inlineStr" in the original OOXML), so force a text number-format before
# writing, then restore the default "Normal" style afterwards so no stray
# formatting is left behind on the cell.
$numericLooking = @('D5', 'D6', 'D8', 'D10', 'D13', 'D18', 'D20', 'D21', 'D23', 'D25', 'D28', 'D29', 'D37', 'D40', 'D42', 'D44', 'D46', 'D49')
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D5').Value = '226.23'
$ws.Range('D6').Value = '0.549'
$ws.Range('D8').Value = '32.31'
$ws.Range('D10').Value = '0.0687'
$ws.Range('D13').Value = '11.01'
$ws.Range('D18').Value = '67.87'
$ws.Range('D20').Value = '245.49'
$ws.Range('D21').Value = '11.02'
$ws.Range('D23').Value = '4.15'
$ws.Range('D25').Value = '161.62'
$ws.Range('D28').Value = '0.115'
$ws.Range('D29').Value = '1.01'
$ws.Range('D37').Value = '0.658'
$ws.Range('D40').Value = '81.76'
$ws.Range('D42').Value = '13.87'
$ws.Range('D44').Value = '0.919'
$ws.Range('D46').Value = '6.11'
$ws.Range('D49').Value = '104.96'

foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining cells: plain text / non-numeric-looking strings, safe to assign directly.
$ws.Range('D2').Value = '34.177.40'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.786.25'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('D12').Value = '2.043.56'
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E13').Value = '  -2.87%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.760.57'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('D16').Value = '34.159.95'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0797'
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('E23').Value = '  +1.02%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E28').Value = '  +1.16%  '
$ws.Range('E29').Value = '  +0.42%  '
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('E32').Value = '  +2.76%  '
$ws.Range('E34').Value = '  -1.82%  '
$ws.Range('D35').Value = '1.443.39'
$ws.Range('E35').Value = '  +2.65%  '
$ws.Range('E36').Value = '  +10.25%  '
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('E41').Value = '  +1.54%  '
$ws.Range('E42').Value = '  +3.69%  '
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('D48').Value = '1.943.60'
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = '0.0₆0129'
$ws.Range('E51').Value = '  -6.40%  '
